$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Change 1: "What are the Locations of the business on campus?" ---
# becomes a spell-checked "vWhat are the Locations of the business on campus?"
$p1 = $d.Paragraphs(1)
$xml1 = @"
<w:p $wns>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:spacing w:line="480" w:lineRule="auto"/>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>v</w:t>
  </w:r>
  <w:r>
    <w:t>What</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> are the Locations of the business on campus?</w:t>
  </w:r>
</w:p>
"@
[void]$p1.Range.InsertXML($xml1)

# --- Change 2: "What does the Random Number mean?" paragraph ---
# becomes "What is cross listed course?"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*does the Random Number mean?*") {
        $target = $p
        break
    }
}

$xml2 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:spacing w:line="480" w:lineRule="auto"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">What is </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">cross </w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>listed course?</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:spacing w:line="480" w:lineRule="auto"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Can we get campus population (in person enrollment) data back to 2015 to match your data?</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:spacing w:line="480" w:lineRule="auto"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Can we get transaction data for all companies (analysis of avg or median price vs number of items to student population and business metrics)</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:spacing w:line="480" w:lineRule="auto"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Can we tell if the business </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:highlight w:val="red"/>
    </w:rPr>
    <w:t>was outside</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">, daytime or </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>evening</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:spacing w:line="480" w:lineRule="auto"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Is there any way to link company letters to the name of the company in the registration data?</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:spacing w:line="480" w:lineRule="auto"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Do the N/</w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>a&#8217;s</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> for company letters and learning teams, faculty means the student dropped the course?   And can we tell when they dropped?</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:spacing w:line="480" w:lineRule="auto"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">What happened to the missing data on the </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Historical Financial Results</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> data?  And are there some semesters with only on IBC company?</w:t>
  </w:r>
</w:p>
"@
[void]$target.Range.InsertXML($xml2)
